$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data table: bus number + three phase voltage readings (A-N, B-N, C-N)
$data = @(
    @(632, "1.0210 at -2.49", "1.0420 at -121.72", "1.0174 at 117.83"),
    @(645, "0.0000 at 0.00",  "1.0371 at -122.38", "1.0101 at 117.40"),
    @(646, "0.0000 at 0.00",  "1.0379 at -122.74", "1.0049 at 117.17"),
    @(633, "1.0180 at -2.55", "1.0401 at -121.77", "1.0148 at 117.83"),
    @(634, "0.9940 at -3.23", "1.0218 at -122.22", "0.9960 at 117.35"),
    @(671, "1.0105 at -4.36", "1.0558 at -121.54", "0.9947 at 117.50"),
    @(684, "1.0103 at -4.40", "0.0000 at 0.00",    "0.9949 at 117.48"),
    @(611, "0.0000 at 0.00",  "0.0000 at 0.00",    "0.9953 at 117.45"),
    @(675, "1.0033 at -4.54", "1.0572 at -121.65", "0.9918 at 117.56"),
    @(680, "1.0105 at -4.36", "1.0558 at -121.54", "0.9947 at 117.50"),
    @(652, "1.0091 at -4.49", "0.0000 at 0.00",    "0.0000 at 0.00")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row = $row + 1
}
